# Applies the author's edit:
#   1. Bumps the cached "datetimeFigureOut" date field text from
#      08-08-2025 to 09-08-2025 everywhere it is cached (the slide
#      master's Date Placeholder and every slide layout's Date
#      Placeholder).
#   2. Fixes the typo in the slide 7 title: CONCLUTION -> CONCLUSION.

$p = $ppt.ActivePresentation

$oldDate = "08-08-2025"
$newDate = "09-08-2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 7 title: fix the "CONCLUTION" typo.
$slide7 = $p.Slides.Item(7)
$title = $slide7.Shapes.Title
if ($title.TextFrame.TextRange.Text -eq "CONCLUTION") {
    $title.TextFrame.TextRange.Text = "CONCLUSION"
}
